# Intervention content and delivery.xlsx
# 1) Insert a new row for "pharmacological treatment" (GMHO:0000262) above the
#    existing row 47 ("physical exertion expended on a behaviour"), shifting
#    every row from the old 47 down through the old 67 down by one.
# 2) Remove the stray empty placeholder cells on row 39 that carried no value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert new row 47 -------------------------------------------------
$ws.Rows.Item(47).Insert()

# Copy formatting (fill style "s=4") from the existing "mental health
# intervention content" row (row 40) onto the freshly inserted row so the new
# row matches the same visual treatment as its sibling rows.
$ws.Range("A40:V40").Copy()
$ws.Range("A47:V47").PasteSpecial(-4122)

$ws.Range("A47").Value = "GMHO:0000262"
$ws.Range("B47").Value = "pharmacological treatment"
$ws.Range("C47").Value = "Mental health intervention content that uses pharmacological substances to assess and improve a person" + [char]0x2019 + "s adaptive mental or behavioural functioning."
$ws.Range("D47").Value = "mental health intervention content"
$ws.Range("J47").Value = "Intervention content and delivery"
$ws.Range("P47").Value = "LSR 1, LSR3"
$ws.Range("S47").Value = "Proposed"

# --- 2. Clean up row 39 ----------------------------------------------------
# These cells previously held no data (empty placeholder cells); remove them
# entirely rather than leaving blank values behind.
$ws.Range("E39:O39").ClearContents()
$ws.Range("R39").ClearContents()
$ws.Range("T39:V39").ClearContents()
